# Apply the "Updated cryptos list" data refresh (GitHub Actions bot commit)
# to Sheet1 of the workbook. Every cell below is a text cell (Price /
# Volume(1h) columns, plus two swapped coin rows), so values are written
# with a leading apostrophe to force Excel to store them as text instead
# of auto-converting number-looking strings (e.g. "1.00", "0.573") into
# numeric values - matching the original inlineStr text cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.860.82"
$ws.Range("E2").Value = "'  +1.11%  "

$ws.Range("D3").Value = "'2.300.43"
$ws.Range("E3").Value = "'  -0.74%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  -0.02%  "

$ws.Range("D5").Value = "'541.64"

$ws.Range("D6").Value = "'129.27"
$ws.Range("E6").Value = "'  -2.53%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  -0.04%  "

$ws.Range("D8").Value = "'0.573"
$ws.Range("E8").Value = "'  -2.40%  "

$ws.Range("D9").Value = "'2.298.73"
$ws.Range("E9").Value = "'  -0.62%  "

$ws.Range("E10").Value = "'  -0.53%  "

$ws.Range("E11").Value = "'  +0.60%  "

$ws.Range("E12").Value = "'  -0.21%  "

$ws.Range("E13").Value = "'  -0.89%  "

$ws.Range("D14").Value = "'23.26"
$ws.Range("E14").Value = "'  -2.93%  "

$ws.Range("D15").Value = "'59.838.38"
$ws.Range("E15").Value = "'  +1.29%  "

$ws.Range("D16").Value = "'2.709.47"
$ws.Range("E16").Value = "'  -0.75%  "

$ws.Range("E17").Value = "'  -1.31%  "

$ws.Range("D18").Value = "'2.294.22"
$ws.Range("E18").Value = "'  -2.01%  "

$ws.Range("D19").Value = "'10.45"
$ws.Range("E19").Value = "'  -1.87%  "

$ws.Range("E20").Value = "'  -2.63%  "

$ws.Range("D21").Value = "'311.05"
$ws.Range("E21").Value = "'  -0.65%  "

$ws.Range("E22").Value = "'  -0.74%  "

$ws.Range("E23").Value = "'  -0.09%  "

$ws.Range("E24").Value = "'  +0.03%  "

$ws.Range("D25").Value = "'63.68"

$ws.Range("D26").Value = "'0.170"
$ws.Range("E26").Value = "'  -1.95%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "'  -0.18%  "

$ws.Range("D28").Value = "'7.71"

$ws.Range("E29").Value = "'  +2.52%  "

$ws.Range("B30").Value = "'SuiNetwork"
$ws.Range("C30").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").Value = "'1.17"
$ws.Range("E30").Value = "'  -0.20%  "

$ws.Range("B31").Value = "'Monero"
$ws.Range("C31").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'170.39"
$ws.Range("E31").Value = "'  -0.01%  "

$ws.Range("E32").Value = "'  -2.22%  "

$ws.Range("E33").Value = "'  -2.59%  "

$ws.Range("D34").Value = "'5.80"
$ws.Range("E34").Value = "'  -1.55%  "

$ws.Range("E35").Value = "'  +1.75%  "

$ws.Range("D36").Value = "'0.378"
$ws.Range("E36").Value = "'  -1.90%  "

$ws.Range("E38").Value = "'  -1.39%  "

$ws.Range("E39").Value = "'  +0.08%  "

$ws.Range("E40").Value = "'  -2.78%  "

$ws.Range("D41").Value = "'317.21"
$ws.Range("E41").Value = "'  +4.18%  "

$ws.Range("D42").Value = "'37.83"
$ws.Range("E42").Value = "'  -1.64%  "

$ws.Range("E43").Value = "'  -1.51%  "

$ws.Range("D44").Value = "'135.69"
$ws.Range("E44").Value = "'  -3.88%  "

$ws.Range("D45").Value = "'3.42"
$ws.Range("E45").Value = "'  -1.25%  "

$ws.Range("E46").Value = "'  -2.49%  "

$ws.Range("E47").Value = "'  +0.67%  "

$ws.Range("D48").Value = "'18.71"
$ws.Range("E48").Value = "'  +1.06%  "

$ws.Range("D49").Value = "'0.0489"
$ws.Range("E49").Value = "'  -1.46%  "

$ws.Range("D50").Value = "'0.0₆0221"
$ws.Range("E50").Value = "'  +19.57%  "

$ws.Range("E51").Value = "'  -0.67%  "
